# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 3.811642989160245 }
    3  = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    4  = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 8.418600821238126 }
    5  = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 8.418600821238126 }
    6  = @{ B = 0.02258322285507441;  C = 0.3375848360084654;  D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.572787616952587 }
    7  = @{ B = 0.1554434735375247;   C = 0.05231270169004087; D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.420375733316612 }
    8  = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    9  = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    10 = @{ B = 0.02258322285507441;  C = 0.3375848360084654;  D = 16.98373111632243;   E = 0.4998867070740569; G = 17.84378588226003 }
    11 = @{ B = 0.1554434735375247;   C = 0.3375848360084654;  D = 0.1529057820181812;  E = 6.48142807727062;    G = 7.127362168834791 }
    12 = @{ B = 3.182878228561681;    C = 87981.0709163148;    D = 0.7127328510149897;  E = 6.48142807727062;    G = 87991.44795547164 }
    13 = @{ B = 0.006876353814593728; C = 0.0001537489499301437; D = 0.7127328510149897; E = 0.4998867070740569; G = 1.219649660853571 }
    14 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 4.371470058157054 }
    15 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 4.371470058157054 }
    16 = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 6.048734245549538 }
    17 = @{ B = 3.182878228561681;    C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 8.418600821238126 }
    18 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 6.741336633845642 }
    19 = @{ B = 1.505614041169197;    C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 4.371470058157054 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $vals.G
}
